$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 2855
    8  = 1533
    12 = 1271
    14 = 401
    16 = 58
    18 = 114
    21 = 2845
    22 = 341
    23 = 10
    24 = 64
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
